$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (date in name moved from 2023-09-08 to 2023-09-09)
$ws.Name = "Product-2023-09-09"

# Update the date/time related string values (createdDate, updatedDate, zonedDateTime columns -> H, I, J)
# Row 2 and Row 3 share the same values for these three columns.
$ws.Range("H2").Value = "09/09/2023 17:45"
$ws.Range("I2").Value = "09 sept. 2023"
$ws.Range("J2").Value = "2023-09-09T17:45:19.6415433+01:00[Africa/Casablanca]"

$ws.Range("H3").Value = "09/09/2023 17:45"
$ws.Range("I3").Value = "09 sept. 2023"
$ws.Range("J3").Value = "2023-09-09T17:45:19.6415433+01:00[Africa/Casablanca]"
